# Update NATMI LR-pair (Cxcl10-Cxcr3) sheet with refreshed TPM-derived
# values: ligand/receptor expression, derived specificities, and edge
# weights, per the "update scripts wuth new tpm" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 21.678587
$ws.Range("H2").Value = 65.03576100000001
$ws.Range("I2").Value = 0.0772399443186744
$ws.Range("J2").Value = 0.07723994431867441
$ws.Range("O2").Value = 0.01611173663836548
$ws.Range("P2").Value = 0.01611173663836548
$ws.Range("Q2").Value = 1.380947670487
$ws.Range("R2").Value = 12.428529034383
$ws.Range("S2").Value = 0.001244469640824496
$ws.Range("T2").Value = 0.001244469640824496
$ws.Range("G3").Value = 21.678587
$ws.Range("H3").Value = 65.03576100000001
$ws.Range("I3").Value = 0.0772399443186744
$ws.Range("J3").Value = 0.07723994431867441
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.701496333333333
$ws.Range("N3").Value = 8.104489000000001
$ws.Range("O3").Value = 0.68328279700753
$ws.Range("P3").Value = 0.68328279700753
$ws.Range("Q3").Value = 58.56462329234768
$ws.Range("R3").Value = 527.0816096311291
$ws.Range("S3").Value = 0.05277672519476972
$ws.Range("T3").Value = 0.05277672519476973
$ws.Range("G4").Value = 21.678587
$ws.Range("H4").Value = 65.03576100000001
$ws.Range("I4").Value = 0.0772399443186744
$ws.Range("J4").Value = 0.07723994431867441
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.188504333333333
$ws.Range("N4").Value = 3.565513
$ws.Range("O4").Value = 0.3006054663541045
$ws.Range("P4").Value = 0.3006054663541044
$ws.Range("Q4").Value = 25.76509459004367
$ws.Range("R4").Value = 231.885851310393
$ws.Range("S4").Value = 0.02321874948308018
$ws.Range("T4").Value = 0.02321874948308018
$ws.Range("I5").Value = 0.4677505770609061
$ws.Range("J5").Value = 0.4677505770609062
$ws.Range("O5").Value = 0.01611173663836548
$ws.Range("P5").Value = 0.01611173663836548
$ws.Range("S5").Value = 0.007536274110048796
$ws.Range("T5").Value = 0.007536274110048797
$ws.Range("I6").Value = 0.4677505770609061
$ws.Range("J6").Value = 0.4677505770609062
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.701496333333333
$ws.Range("N6").Value = 8.104489000000001
$ws.Range("O6").Value = 0.68328279700753
$ws.Range("P6").Value = 0.68328279700753
$ws.Range("Q6").Value = 354.656345003698
$ws.Range("R6").Value = 3191.907105033282
$ws.Range("S6").Value = 0.3196059225960621
$ws.Range("T6").Value = 0.3196059225960622
$ws.Range("I7").Value = 0.4677505770609061
$ws.Range("J7").Value = 0.4677505770609062
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.188504333333333
$ws.Range("N7").Value = 3.565513
$ws.Range("O7").Value = 0.3006054663541045
$ws.Range("P7").Value = 0.3006054663541044
$ws.Range("Q7").Value = 156.028567457266
$ws.Range("R7").Value = 1404.257107115394
$ws.Range("S7").Value = 0.1406083803547951
$ws.Range("T7").Value = 0.1406083803547952
$ws.Range("G8").Value = 67.370907
$ws.Range("H8").Value = 202.112721
$ws.Range("I8").Value = 0.2400398653924534
$ws.Range("J8").Value = 0.2400398653924535
$ws.Range("O8").Value = 0.01611173663836548
$ws.Range("P8").Value = 0.01611173663836548
$ws.Range("Q8").Value = 4.291594146806999
$ws.Range("R8").Value = 38.624347321263
$ws.Range("S8").Value = 0.003867459093911909
$ws.Range("T8").Value = 0.00386745909391191
$ws.Range("G9").Value = 67.370907
$ws.Range("H9").Value = 202.112721
$ws.Range("I9").Value = 0.2400398653924534
$ws.Range("J9").Value = 0.2400398653924535
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.701496333333333
$ws.Range("N9").Value = 8.104489000000001
$ws.Range("O9").Value = 0.68328279700753
$ws.Range("P9").Value = 0.68328279700753
$ws.Range("Q9").Value = 182.002258233841
$ws.Range("R9").Value = 1638.020324104569
$ws.Range("S9").Value = 0.1640151106186666
$ws.Range("T9").Value = 0.1640151106186666
$ws.Range("G10").Value = 67.370907
$ws.Range("H10").Value = 202.112721
$ws.Range("I10").Value = 0.2400398653924534
$ws.Range("J10").Value = 0.2400398653924535
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.188504333333333
$ws.Range("N10").Value = 3.565513
$ws.Range("O10").Value = 0.3006054663541045
$ws.Range("P10").Value = 0.3006054663541044
$ws.Range("Q10").Value = 80.07061491009701
$ws.Range("R10").Value = 720.6355341908732
$ws.Range("S10").Value = 0.07215729567987492
$ws.Range("T10").Value = 0.07215729567987492
$ws.Range("G11").Value = 8.783890333333334
$ws.Range("H11").Value = 26.351671
$ws.Range("I11").Value = 0.03129665232554173
$ws.Range("J11").Value = 0.03129665232554174
$ws.Range("O11").Value = 0.01611173663836548
$ws.Range("P11").Value = 0.01611173663836548
$ws.Range("Q11").Value = 0.5595425981236667
$ws.Range("R11").Value = 5.035883383113
$ws.Range("S11").Value = 0.000504243419931617
$ws.Range("T11").Value = 0.000504243419931617
$ws.Range("G12").Value = 8.783890333333334
$ws.Range("H12").Value = 26.351671
$ws.Range("I12").Value = 0.03129665232554173
$ws.Range("J12").Value = 0.03129665232554174
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.701496333333333
$ws.Range("N12").Value = 8.104489000000001
$ws.Range("O12").Value = 0.68328279700753
$ws.Range("P12").Value = 0.68328279700753
$ws.Range("Q12").Value = 23.72964752790211
$ws.Range("R12").Value = 213.566827751119
$ws.Range("S12").Value = 0.02138446413796838
$ws.Range("T12").Value = 0.02138446413796838
$ws.Range("G13").Value = 8.783890333333334
$ws.Range("H13").Value = 26.351671
$ws.Range("I13").Value = 0.03129665232554173
$ws.Range("J13").Value = 0.03129665232554174
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.188504333333333
$ws.Range("N13").Value = 3.565513
$ws.Range("O13").Value = 0.3006054663541045
$ws.Range("P13").Value = 0.3006054663541044
$ws.Range("Q13").Value = 10.43969172469145
$ws.Range("R13").Value = 93.957225522223
$ws.Range("S13").Value = 0.00940794476764174
$ws.Range("T13").Value = 0.009407944767641742
$ws.Range("G14").Value = 51.550662
$ws.Range("H14").Value = 154.651986
$ws.Range("I14").Value = 0.1836729609024243
$ws.Range("J14").Value = 0.1836729609024243
$ws.Range("O14").Value = 0.01611173663836548
$ws.Range("P14").Value = 0.01611173663836548
$ws.Range("Q14").Value = 3.283828720062
$ws.Range("R14").Value = 29.554458480558
$ws.Range("S14").Value = 0.00295929037364866
$ws.Range("T14").Value = 0.002959290373648659
$ws.Range("G15").Value = 51.550662
$ws.Range("H15").Value = 154.651986
$ws.Range("I15").Value = 0.1836729609024243
$ws.Range("J15").Value = 0.1836729609024243
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.701496333333333
$ws.Range("N15").Value = 8.104489000000001
$ws.Range("O15").Value = 0.68328279700753
$ws.Range("P15").Value = 0.68328279700753
$ws.Range("Q15").Value = 139.263924373906
$ws.Range("R15").Value = 1253.375319365154
$ws.Range("S15").Value = 0.1255005744600632
$ws.Range("T15").Value = 0.1255005744600632
$ws.Range("G16").Value = 51.550662
$ws.Range("H16").Value = 154.651986
$ws.Range("I16").Value = 0.1836729609024243
$ws.Range("J16").Value = 0.1836729609024243
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.188504333333333
$ws.Range("N16").Value = 3.565513
$ws.Range("O16").Value = 0.3006054663541045
$ws.Range("P16").Value = 0.3006054663541044
$ws.Range("Q16").Value = 61.268185173202
$ws.Range("R16").Value = 551.413666558818
$ws.Range("S16").Value = 0.05521309606871244
$ws.Range("T16").Value = 0.05521309606871244
